$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 to the "新题" sheet: date, problem title ("144 tree" -
# LeetCode 144 Binary Tree Preorder Traversal, per the commit message),
# and a "done" mark in the 帅 (E) column.
#
# Copy A14's formatting (the existing short-date number style) down to A15
# first, so the new date cell reuses the same style instead of Excel
# inferring a brand-new number format for it.
$ws.Range("A14").Copy($ws.Range("A15"))
$ws.Range("A15").Value = "3/24/2019"
$ws.Range("B15").Value = "144 tree"
$ws.Range("E15").Value = "done"

# Match the author's final selection/active cell.
$ws.Range("B17").Select()
